# Insert a new weekly price-report row at row 47 (pushing existing rows 47-123
# down to 48-124) and populate it with the new "Haba" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 47; everything below (old rows 47-123)
# shifts down to 48-124, carrying its formatting with it.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data point. The
# non-changing columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the same template
# values used throughout the rest of this "Haba" sheet.
$ws.Range("A47").Value = 2
$ws.Range("B47").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 45246
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 100112026
$ws.Range("G47").Value = "Haba"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 600
$ws.Range("K47").Value = 8000
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = 8500
$ws.Range("N47").Value = "`$/saco 25 kilos"
$ws.Range("O47").Value = "Provincia de Limarí"
$ws.Range("P47").Value = 340
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

# Match the date style used by the rest of column D.
$ws.Range("D47").NumberFormat = $ws.Range("D48").NumberFormat
